# Update the metrics table (rows 2-26, columns B:Q) with the new values
# produced after retraining with the new LM (per commit message
# "atualizado todo o treinamento para o novo lm").
# Every data row previously shared identical values across B:Q, and after
# the update they again share identical (but different) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    [double]"0.9999677858346675",
    [double]"0.9990493920370204",
    [double]"0.999994540723036",
    [double]"0.9999019866292851",
    [double]"0.9999286684456561",
    [double]"3.007050929946918e-05",
    [double]"0.0008873508065744985",
    [double]"2.739584952253493e-06",
    [double]"0.0001361276375442473",
    [double]"6.943361116717897e-05",
    [double]"0.0003164627591308915",
    [double]"0.005483658386466938",
    [double]"1.000026659998896",
    [double]"0.005717109076716589",
    [double]"126.8239312479565",
    [double]"191.4243499659711"
)

$firstRow = 2
$lastRow = 26
$firstCol = 2  # Column B
$lastCol = 17  # Column Q

for ($row = $firstRow; $row -le $lastRow; $row++) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - $firstCol]
    }
}
